$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text representation
# instead of being auto-converted to a number by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.484.04"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.727.94"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.44"
$ws.Range("E5").Value = "  +2.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4796"
$ws.Range("E7").Value = "  +1.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2691"
$ws.Range("E8").Value = "  +1.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06224"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.729.53"
$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07136"
$ws.Range("E11").Value = "  +0.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.74"
$ws.Range("E12").Value = "  +3.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6191"
$ws.Range("E13").Value = "  +4.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.542"
$ws.Range("E14").Value = "  +2.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.26"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.482.51"
$ws.Range("E17").Value = "  +0.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9988"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006952"
$ws.Range("E19").Value = "  +2.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.71"
$ws.Range("E20").Value = "  +0.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.952.83"
$ws.Range("E21").Value = "  +1.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.542"
$ws.Range("E22").Value = "  -0.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.936"
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.304"
$ws.Range("E24").Value = "  -0.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.38"
$ws.Range("E25").Value = "  +0.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.41"
$ws.Range("E26").Value = "  +1.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.804"
$ws.Range("E27").Value = "  +2.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.405"
$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.05"
$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.984"
$ws.Range("E30").Value = "  -0.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08036"
$ws.Range("E31").Value = "  +3.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.733"
$ws.Range("E32").Value = "  +1.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04562"
$ws.Range("E33").Value = "  +3.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.617"
$ws.Range("E34").Value = "  +0.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6387"
$ws.Range("E35").Value = "  +2.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9907"
$ws.Range("E36").Value = "  +2.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9340"
$ws.Range("E37").Value = "  +1.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.056"
$ws.Range("E38").Value = "  +7.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.416"
$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.21"
$ws.Range("E40").Value = "  -4.79%  "

$ws.Range("E41").Value = "  +0.22%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.719"
$ws.Range("E42").Value = "  +10.93%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01502"
$ws.Range("E43").Value = "  +2.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3919"
$ws.Range("E44").Value = "  +2.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.996"
$ws.Range("E45").Value = "  +12.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1192"
$ws.Range("E46").Value = "  +4.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05319"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "31.11"
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.929"
$ws.Range("E49").Value = "  +3.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.269"
$ws.Range("E50").Value = "  +3.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3428"
$ws.Range("E51").Value = "  +1.40%  "
